# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G ("K") holds the strikeout count for each start. The source data
# pull was regenerated, so every row's K value is recalculated here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values, in row order starting at row 2 (most recent start) through row 27.
$kValues = @(4, 2, 3, 4, 5, 5, 3, 3, 1, 7, 6, 8, 4, 9, 5, 3, 4, 3, 2, 2, 2, 2, 3, 5, 5, 1)

$row = 2
foreach ($k in $kValues) {
    $ws.Cells.Item($row, 7).Value = $k
    $row++
}
